$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.08029043674469
$ws.Range("B1").Value = 1.445855379104614
$ws.Range("C1").Value = 2.355828762054443
$ws.Range("D1").Value = 4.56023645401001
$ws.Range("E1").Value = 1.934415459632874
